# List of bugs updated, few minor bugs fixed.
# Adds three new bug rows (8, 9, 10) to the WinArc bug-report sheet and
# marks the last of them as Fixed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: bug #8 -----------------------------------------------------
$ws.Range("A16").Value = 8
$ws.Range("A16").HorizontalAlignment = -4108   # xlCenter - matches ID column style
$ws.Range("B16").Value = "High"
$ws.Range("C16").Value = "Creating zip-archive, that contains folder."
$ws.Range("C16").VerticalAlignment = -4108     # xlCenter (vertical) - matches DESCRIPTION column style

# --- Row 17: bug #9 -------------------------------------------------------
$ws.Range("A17").Value = 9
$ws.Range("A17").HorizontalAlignment = -4108
$ws.Range("B17").Value = "Low"
$ws.Range("C17").Value = "Folder icons in file system view are absent."
$ws.Range("C17").VerticalAlignment = -4108

# --- Row 18: bug #10 (fixed) ----------------------------------------------
$ws.Range("A18").Value = 10
$ws.Range("A18").HorizontalAlignment = -4108
$ws.Range("B18").Value = "Low"
$ws.Range("C18").Value = "There isn't any constant local that contains number of columns for folder view."
$ws.Range("C18").VerticalAlignment = -4108
$ws.Range("D18").Value = "Fixed"

# New STATUS column width for column D (matches width added next to the table)
$ws.Columns.Item(4).ColumnWidth = 12.8

# Update selection / active cell to the new last entry, like Excel would
# leave it after typing the last value.
$ws.Range("D18").Select() | Out-Null
